$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.766.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.766.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.546"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.020.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.790.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.729.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.612"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0774"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0508"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.379.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.652"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0185"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.903"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.17%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0137"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.922.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
